$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '21.290.95'
$ws.Range('E2').Value = '  -3.29%  '

# Row 3
$ws.Range('D3').Value = '1.522.36'
$ws.Range('E3').Value = '  -1.92%  '

# Row 4
$ws.Range('E4').Value = '  +0.28%  '

# Row 5
$ws.Range('D5').Value = '''1.005'
$ws.Range('E5').Value = '  +0.32%  '

# Row 6
$ws.Range('D6').Value = '''288.22'
$ws.Range('E6').Value = '  -0.69%  '

# Row 7
$ws.Range('D7').Value = '''0.3895'
$ws.Range('E7').Value = '  -0.81%  '

# Row 8
$ws.Range('D8').Value = '''0.3157'
$ws.Range('E8').Value = '  -1.61%  '

# Row 9
$ws.Range('D9').Value = '''42.55'
$ws.Range('E9').Value = '  -2.06%  '

# Row 10
$ws.Range('D10').Value = '''0.07071'
$ws.Range('E10').Value = '  -2.31%  '

# Row 11
$ws.Range('D11').Value = '''1.062'
$ws.Range('E11').Value = '  -0.67%  '

# Row 12
$ws.Range('D12').Value = '''1.007'
$ws.Range('E12').Value = '  +0.41%  '

# Row 13
$ws.Range('D13').Value = '''5.670'
$ws.Range('E13').Value = '  +0.23%  '

# Row 14
$ws.Range('D14').Value = '''18.03'
$ws.Range('E14').Value = '  -3.29%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.528.80'
$ws.Range('E15').Value = '  -1.48%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '''6.428'
$ws.Range('E16').Value = '  -2.58%  '

# Row 17
$ws.Range('D17').Value = '''0.00001082'
$ws.Range('E17').Value = '  -3.87%  '

# Row 18
$ws.Range('D18').Value = '''0.06608'
$ws.Range('E18').Value = '  +0.37%  '

# Row 19
$ws.Range('D19').Value = '''82.53'
$ws.Range('E19').Value = '  -0.89%  '

# Row 20
$ws.Range('D20').Value = '''1.003'
$ws.Range('E20').Value = '  +0.21%  '

# Row 21
$ws.Range('D21').Value = '''6.065'
$ws.Range('E21').Value = '  -3.25%  '

# Row 22
$ws.Range('D22').Value = '''15.27'
$ws.Range('E22').Value = '  -0.97%  '

# Row 23
$ws.Range('D23').Value = '''10.77'
$ws.Range('E23').Value = '  -4.04%  '

# Row 24
$ws.Range('D24').Value = '''2.367'
$ws.Range('E24').Value = '  -0.51%  '

# Row 25
$ws.Range('D25').Value = '21.273.62'
$ws.Range('E25').Value = '  -3.43%  '

# Row 26
$ws.Range('D26').Value = '''2.374'
$ws.Range('E26').Value = '  -1.49%  '

# Row 27
$ws.Range('D27').Value = '''147.69'
$ws.Range('E27').Value = '  -0.68%  '

# Row 28
$ws.Range('D28').Value = '''18.26'
$ws.Range('E28').Value = '  -1.31%  '

# Row 29
$ws.Range('D29').Value = '''4.815'
$ws.Range('E29').Value = '  -1.38%  '

# Row 30
$ws.Range('D30').Value = '1.697.73'
$ws.Range('E30').Value = '  -1.70%  '

# Row 31
$ws.Range('D31').Value = '''116.18'

# Row 32
$ws.Range('D32').Value = '''6.017'
$ws.Range('E32').Value = '  +3.86%  '

# Row 33
$ws.Range('D33').Value = '''0.9544'
$ws.Range('E33').Value = '  -3.61%  '

# Row 34
$ws.Range('D34').Value = '''0.08010'
$ws.Range('E34').Value = '  -3.30%  '

# Row 35
$ws.Range('D35').Value = '''8.508'
$ws.Range('E35').Value = '  -5.27%  '

# Row 36
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.496'
$ws.Range('E36').Value = '  -7.40%  '

# Row 37
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '''5.119'
$ws.Range('E37').Value = '  +0.59%  '

# Row 38
$ws.Range('D38').Value = '''11.40'
$ws.Range('E38').Value = '  +7.25%  '

# Row 39
$ws.Range('D39').Value = '''0.05925'
$ws.Range('E39').Value = '  -2.28%  '

# Row 40
$ws.Range('D40').Value = '''0.02161'
$ws.Range('E40').Value = '  -4.10%  '

# Row 41
$ws.Range('D41').Value = '''0.1999'
$ws.Range('E41').Value = '  -1.79%  '

# Row 42
$ws.Range('D42').Value = '''1.171'
$ws.Range('E42').Value = '  -3.41%  '

# Row 43
$ws.Range('D43').Value = '''1.004'
$ws.Range('E43').Value = '  +0.45%  '

# Row 44
$ws.Range('D44').Value = '''0.5698'
$ws.Range('E44').Value = '  -1.58%  '

# Row 45
$ws.Range('D45').Value = '''12.98'
$ws.Range('E45').Value = '  +0.54%  '

# Row 46
$ws.Range('D46').Value = '''3.714'
$ws.Range('E46').Value = '  -0.86%  '

# Row 47
$ws.Range('D47').Value = '''0.5516'
$ws.Range('E47').Value = '  -0.79%  '

# Row 48
$ws.Range('D48').Value = '''1.880'
$ws.Range('E48').Value = '  -0.47%  '

# Row 49
$ws.Range('D49').Value = '''1.148'
$ws.Range('E49').Value = '  +1.48%  '

# Row 50
$ws.Range('D50').Value = '''115.33'
$ws.Range('E50').Value = '  -2.35%  '

# Row 51
$ws.Range('D51').Value = '''0.06598'
$ws.Range('E51').Value = '  -3.15%  '
